$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Location"
$ws.Range("C1").Style = "Normal"

$ws.Range("A12").Value = "Pune"
$ws.Range("A12").Style = "Normal"

$ws.Range("O6").Select()
